$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.424211666666666
$ws.Range("H2").Value = 4.272634999999999
$ws.Range("M2").Value = 0.2848286666666667
$ws.Range("N2").Value = 0.854486
$ws.Range("O2").Value = 0.08022967564521397
$ws.Range("P2").Value = 0.08022967564521397
$ws.Range("Q2").Value = 0.4056563100677777
$ws.Range("R2").Value = 3.650906790609999
$ws.Range("S2").Value = 0.08022967564521397
$ws.Range("T2").Value = 0.08022967564521397

$ws.Range("G3").Value = 1.424211666666666
$ws.Range("H3").Value = 4.272634999999999
$ws.Range("O3").Value = 0.2730594381596053
$ws.Range("P3").Value = 0.2730594381596053
$ws.Range("Q3").Value = 1.380639809673889
$ws.Range("R3").Value = 12.425758287065
$ws.Range("S3").Value = 0.2730594381596053
$ws.Range("T3").Value = 0.2730594381596053

$ws.Range("G4").Value = 1.424211666666666
$ws.Range("H4").Value = 4.272634999999999
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.351483666666667
$ws.Range("N4").Value = 4.054451
$ws.Range("O4").Value = 0.3806818235166093
$ws.Range("P4").Value = 0.3806818235166093
$ws.Range("Q4").Value = 1.924798805376111
$ws.Range("R4").Value = 17.323189248385
$ws.Range("S4").Value = 0.3806818235166093
$ws.Range("T4").Value = 0.3806818235166093

$ws.Range("G5").Value = 1.424211666666666
$ws.Range("H5").Value = 4.272634999999999
$ws.Range("M5").Value = 0.23571
$ws.Range("N5").Value = 0.70713
$ws.Range("O5").Value = 0.066394078474077
$ws.Range("P5").Value = 0.066394078474077
$ws.Range("Q5").Value = 0.3357009319499999
$ws.Range("R5").Value = 3.02130838755
$ws.Range("S5").Value = 0.066394078474077
$ws.Range("T5").Value = 0.066394078474077

$ws.Range("G6").Value = 1.424211666666666
$ws.Range("H6").Value = 4.272634999999999
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.7087373333333332
$ws.Range("N6").Value = 2.126212
$ws.Range("O6").Value = 0.1996349842044944
$ws.Range("P6").Value = 0.1996349842044945
$ws.Range("Q6").Value = 1.009391978735555
$ws.Range("R6").Value = 9.084527808619997
$ws.Range("S6").Value = 0.1996349842044944
$ws.Range("T6").Value = 0.1996349842044945
